$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fill in the "Informations generales" table (Groupe / Membres / Theme)
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)

$t.Cell(1, 2).Range.Text = "OnlyUpSàrl"
$t.Cell(2, 2).Range.Text = "Eliott – Dioussé - Nicola"
$t.Cell(3, 2).Range.Text = " Immeuble pour habitations et bureaux"

# ------------------------------------------------------------------
# 2. Append the sign-off block at the end of the document:
#      3 empty "Corps de texte" paragraphs, then
#      "Lu et approuvé, 23.01.24, Scherrer Eliott" (also "Corps de texte"),
#      with a _GoBack bookmark right after the inserted text.
#
#    NOTE: use $d.Content.Paragraphs (a fresh Range-derived collection)
#    rather than the cached $d.Paragraphs after the table edits above —
#    the cached collection's indexing gets confused once a table cell's
#    Range.Text has been reassigned.
# ------------------------------------------------------------------
$paras = $d.Content.Paragraphs
$lastPara = $paras.Item($paras.Count)
$r = $lastPara.Range

$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()

$paras2 = $d.Content.Paragraphs
$signPara = $paras2.Item($paras2.Count)
$sr = $signPara.Range
$sr.Collapse(0)
$sr.InsertBefore("Lu et approuvé, ")
$sr.Collapse(0)
# trailing "X" is a throw-away placeholder that lets us anchor the
# _GoBack bookmark exactly after the inserted text (before the pilcrow);
# a zero-length range always snaps to document-start when bookmarked.
$sr.InsertBefore("23.01.24, Scherrer EliottX")

$paras3 = $d.Content.Paragraphs
$signPara2 = $paras3.Item($paras3.Count)
$sr2 = $signPara2.Range
$posEnd = $sr2.End - 1
$placeholderRange = $d.Range($posEnd - 1, $posEnd)
$d.Bookmarks.Add("_GoBack", $placeholderRange)
$placeholderRange2 = $d.Range($posEnd - 1, $posEnd)
$placeholderRange2.Delete()
